# Supplemental Table 5 (phylogeny information) - diversity stats update for Fig 1
#
# Updates several numeric/model values in the big summary table. Each edit
# is applied by directly setting the text of the target table cell's Range
# (rather than Find/Replace), because Find.Execute on this runtime is not
# bounded to the Range it's called on -- it searches/replaces across the
# whole document, which would corrupt unrelated cells that happen to
# contain the same digits as a substring (e.g. "54" inside "2788549").
#
# Row layout (Tables(1)):
#   Row 2  -> "Summary polymerase" (Fig 2A)
#   Row 8  -> "Unclassified bat picornavirus" (Fig 2G)
#   Row 9  -> "Sapelovirus" (Fig 2H)
#   Row 10 -> "Teschovirus" (Fig 2I)
#   Row 11 -> "Sapovirus" (Fig 2J)
#
# Column layout:
#   4 -> # Novel seq
#   5 -> # Reference seq
#   8 -> Overlap length (bp)
#   9 -> Best model

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

function Set-CellText($table, $row, $col, $oldText, $newText) {
    $cell = $table.Cell($row, $col)
    $cellRange = $cell.Range
    # A table cell's Range.Text ends with a paragraph mark + cell mark
    # (two characters) that must be excluded before reassigning .Text.
    $fullText = $cellRange.Text
    $markLen = 2
    $textOnlyLen = $fullText.Length - $markLen
    $textOnly = $fullText.Substring(0, $textOnlyLen)

    if ($textOnly -ne $oldText) {
        throw "Cell ($row,$col) text mismatch: expected '$oldText' but found '$textOnly'"
    }

    $textRange = $d.Range($cellRange.Start, $cellRange.Start + $textOnlyLen)
    $textRange.Text = $newText
}

# Row 2 - Summary polymerase (Fig 2A)
Set-CellText $t 2 5 "267" "273"
Set-CellText $t 2 9 "GTR+G4" "GTR+I+G4"

# Row 8 - Unclassified bat picornavirus (Fig 2G)
Set-CellText $t 8 4 "5" "4"
Set-CellText $t 8 5 "26" "35"
Set-CellText $t 8 8 "~2300" "~1700"

# Row 9 - Sapelovirus (Fig 2H)
Set-CellText $t 9 5 "54" "57"

# Row 10 - Teschovirus (Fig 2I)
Set-CellText $t 10 5 "33" "35"

# Row 11 - Sapovirus (Fig 2J)
Set-CellText $t 11 5 "213" "217"
